# planety_quizz.xlsx - "otazky -> vylepsenie vyhodnotenia"
#
# Replaces the placeholder question text (column B, rows 11-90) that was
# previously re-using the "Otazka slnko N" / "Spravna N" style strings from
# the SLNKO block, with per-planet question text ("Otazka Merkur 1", ...,
# "otazka jupiter 10", "Otazka neptun 10", etc.), and nudges the sheet's
# selection/scroll position to reflect where editing finished.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of the 8 ten-row blocks (rows 11-20, 21-30, ... 81-90) corresponds to
# one planet (A column already holds the planet name); update the question
# text in column B for each of the 10 questions in that block. Casing of the
# "Otazka"/"otazka" prefix and the planet label matches the source data
# (jupiter's block uses an all-lowercase "otazka jupiter N").
$sections = @(
    @{ Prefix = "Otazka"; Label = "Merkur" },
    @{ Prefix = "Otazka"; Label = "Venusa" },
    @{ Prefix = "Otazka"; Label = "Zem" },
    @{ Prefix = "Otazka"; Label = "mars" },
    @{ Prefix = "otazka"; Label = "jupiter" },
    @{ Prefix = "Otazka"; Label = "saturn" },
    @{ Prefix = "Otazka"; Label = "uran" },
    @{ Prefix = "Otazka"; Label = "neptun" }
)

$row = 11
foreach ($sec in $sections) {
    for ($i = 1; $i -le 10; $i++) {
        $ws.Cells.Item($row, 2).Value = "$($sec.Prefix) $($sec.Label) $i"
        $row++
    }
}

# Reflect the final cursor/scroll position left by the edit: scrolled down so
# row 67 is at the top of the view, with B93 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws.Range("B93").Select()
